$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 7) to the experiments table, mirroring the
# structure of the existing rows (timestamp, numeric duration, and a set
# of stringified-list columns).
$ws.Range("A7").Value = '2024-08-11 13:14:47'
$ws.Range("B7").Value = 2.00046181678772
$ws.Range("C7").Value = '[89.25735735899329, 18.266877742822246, 28.559313904069086, 63.23096816831339, 5.9320092705275425, 81.83549357910616, 27.68251191236182, 97.43208525001323, 17.126677661430055, 74.89561560331877]'
$ws.Range("D7").Value = '[38.62247849692741, 37.807462399642624, 5.273182827729178, 34.208400558141534, 34.30295324556304, 33.23412965177267, 22.414082018829806, 35.16621785783049, 9.894010639383072, 36.71634761566201]'
$ws.Range("E7").Value = '[32.981443845428956, 80.11513056673046, 83.8602450408359, 90.76745155683113, 68.15791291305707, 36.27085032026224, 41.86835288966364, 51.308494883311525, 38.74361330891567, 79.58259092638937]'
$ws.Range("F7").Value = '[9.968786189872151, 5.586090523558925, 9.014276506040146, 9.267152408120594, 9.641027572855966, 8.45960925656044, 4.0114419356892315, 5.821995162876144, 8.557301420022632, 2.8964280785815695]'
$ws.Range("G7").Value = '[0.5212153887701326, 0.7040077727468171, 0.5620589154267042, 0.3706984486857198, 0.9775871241710028, 0.26994042708418386, 0.6066561233460374, 0.923203267926135, 0.8528467521626611, 0.3272424483172215]'
$ws.Range("H7").Value = '[511.01848359513065, 468.97432911383817, 831.8901109731171, 397.2120474609366, 877.1744802293973, 234.56235489445947, 912.4792681389679, 668.0175637809969, 695.8846180835076, 379.9392921792524]'
$ws.Range("I7").Value = '[10.542762391010202, 19.142322430162974, 3.5156264612478023, 3.731404382165208, 5.132939960366723, 13.730580785565714, 11.609573166294052, 6.536986144208064, 6.884638435496864, 17.071968493358252]'
$ws.Range("J7").Value = '[8741.941266656155, 6729.149552077638, 7536.365506208702, 2515.7268464027366, 5364.58040405793, 4603.547733687761, 8968.559318745067, 2086.7021741055887, 5555.169518354282, 7163.792291003467]'
$ws.Range("K7").Value = '[1499.6865011555358, 2245.170378198337, 582.9993544455099, 2335.384881360709, 4983.544111246324, 3471.8910598325942, 2654.4523984493935, 2471.0041773311073, 4037.5944853304422, 3888.0678193239473]'
$ws.Range("L7").Value = '[475.8613861664442, 722.0181909427934, 433.0133022632603, 236.25585097905108, 105.10183925862664, 560.538171471306, 229.48404873992135, 158.15037831925528, 766.9390525421402, 106.88711605810329]'
$ws.Range("M7").Value = '[94354.6504388469, 21606.730964200615, 18589.68155822737, 86435.5334093719, 42678.11520929945, 15008.15133250635, 60340.625945523294, 22649.515922699873, 65427.05833271308, 52909.2653954277]'
$ws.Range("N7").Value = '[75.01292321029865, 54.043014819380794, 13.62862312103248, 26.465302230209687, 62.077810323819534, 17.715440645267055, 72.07094684073517, 72.96133620905309, 51.45619525245194, 86.1253889075882]'
$ws.Range("O7").Value = '[1.5997452330569706, 0.5854169093810979, 12.617256905317204, 17.83030852321273, 8.762541726975437, 38.62074351316792, 31.274856942333273, 35.02769163323579, 33.75190280978166, 0.0806925446184048]'
$ws.Range("P7").Value = '[4.71645022754566, 16.473583704029572, 35.15133343257506, 45.74188719875158, 9.67338218096928, 25.14073529103504, 7.773124586758345, 41.43392816849628, 49.22556848352685, 0.3207038058366307]'
$ws.Range("Q7").Value = '[3, 1, 3, 4, 3, 2, 3, 1, 3, 1]'
$ws.Range("R7").Value = '[8.969110363651932, 14.331630428153979, 1.6328413395175079, 6.338491720220105, 9.703050953544567, 9.308227243350052, 10.95489506453281, 11.042933656418281, 12.758479103935166, 13.765271673656187]'
$ws.Range("S7").Value = '[68.56692437564124, 99.05747453183234, 52.802874894626385, 26.45543495983196, 46.157531993622285, 39.56253238819318, 24.03177766369707, 51.236573356632185, 95.1980279745234, 25.814219503605727]'
$ws.Range("T7").Value = '[368.5522213294286, 370.0062980766466, 474.9231040065249, 489.70038441528624, 106.16510652381747, 438.10520118593115, 490.35145972743504, 192.94510060322295, 167.5479743952091, 216.316447494333]'
$ws.Range("U7").Value = '[339.872728722169, 329.9314804578963, 217.06281079300868, 382.03176809429294, 432.7826036226013, 233.1193702883852, 102.22201399950306, 453.70736278052163, 109.2736697469796, 277.3890906947289]'
